# Montana NCES school dictionary: the "flathead h s" row (row 9) mapped to
# the same ACFR entity (Kalispell Public Schools) as "kalispell elem" (row 8),
# so multiple NCES school records map to a single ACFR. The row's data is
# cleared (keeping the highlighted fill style on A9:D9) so it no longer
# duplicates the ACFR mapping, consolidating multiple Montana NCES schools
# under a single ACFR record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the whole row first (mirrors selecting row 9 in the UI before
# clearing it), then clear its contents.
$row9 = $ws.Rows.Item(9)
$row9.Select()
$row9.ClearContents()
